$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# The diagonal cells (B2, C3, D4, ... N14) were highlighted yellow and held a
# literal "1" (the self-correlation of each variable). Remove both: clear the
# value and remove the now-unused yellow highlight fill, while leaving the
# cell's style slot (borders/number format) otherwise intact.
$diagonal = @("B2","C3","D4","E5","F6","G7","H8","I9","J10","K11","L12","M13","N14")

foreach ($addr in $diagonal) {
    $cell = $ws.Range($addr)
    $cell.ClearContents()
    $cell.Interior.ColorIndex = -4142   # xlColorIndexNone
    $cell.Interior.Pattern = -4142      # xlNone
}

# Restore the selection to the cell it was left on when the file was saved.
$ws.Range("N14").Select()
